# MasterTable_process_b_y2.xlsx correction pass.
# - Replace the old ML-model tags (MLA / MLB) in column F with the new,
#   more specific optimizer configuration tags.
# - Move the active selection to K6 (was D13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Column F ("USE_ACTUAL_MODEL") - rows 2-9.
# Old "MLB" rows split into PR_B_Y2/PR_B_Y3 (secundary/primary inputs) and
# PR_B_Y2 (observed inputs); the target row gets PR_B_Y2/PR_C_Y2.
$ws.Range("F2").Value = "PR_B_Y2/PR_B_Y3"
$ws.Range("F3").Value = "PR_B_Y2/PR_B_Y3"
$ws.Range("F4").Value = "PR_B_Y2"
$ws.Range("F5").Value = "PR_B_Y2"
$ws.Range("F6").Value = "PR_B_Y2/PR_C_Y2"

# Old "MLA" rows (W1/W2/W3 features) now also use PR_B_Y2.
$ws.Range("F7").Value = "PR_B_Y2"
$ws.Range("F8").Value = "PR_B_Y2"
$ws.Range("F9").Value = "PR_B_Y2"

# Update the active cell selection shown when the sheet is reopened.
[void]$ws.Range("K6").Select()
